$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "nadpisy v excel kep3" - update the two table headings ("current" / "new" -> "original" / "rewrite1")
$ws.Range("A2").Value = "původní"
$ws.Range("E2").Value = "přepis1"

# match the author's final cursor position in the saved file
$ws.Range("E3").Select()
